$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "https://images4.jdmagicbox.com/uae/jdcatalogue/dubai/77/2013061677/catalogue/a8b1f7ed9ca7c791883ada40a6b0fef9.jpg?output-quality=100"
$ws.Range("E9").Value = "https://images3.jdmagicbox.com/uae/jdcatalogue/abu_dhabi/88/2013123988/catalogue/e18ffdfba70cc9aceff9a07911af28d5.jpg?output-quality=100"
$ws.Range("E11").Value = "https://images3.jdmagicbox.com/uae/jdcatalogue/abu_dhabi/77/nde0244705677/catalogue/3969920607d9a5e15a142f8f74c998a7.jpg?output-quality=100"
